# Loan RBI, Variable Instalments
# - Insert a new (blank) column before column N on the "Repayment Schedule"
#   sheet, shifting the old N/O/P ("Late", "Heading", "Outstanding") columns
#   one position to the right (O/P/Q).
# - Make "Repayment Schedule" the active sheet / update the selection there.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Activating the sheet updates workbook.xml's activeTab and moves
# tabSelected between the sheetViews of "Repayment Schedule" and
# "Transactions".
$ws.Activate()

# Insert a new blank column at N (pushes existing N:P to O:Q).
$ws.Columns("N").Insert()

# Give the freshly inserted column the same width as its neighbour (column M).
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Update the selected cell on the Repayment Schedule sheet.
$ws.Range("J16").Select()
